$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''34.982.98'
$ws.Range('D3').Value = '''1.846.01'
$ws.Range('E3').Value = '  +2.03%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').Value = '''233.03'
$ws.Range('E5').Value = '  +0.66%  '
$ws.Range('E6').Value = '  +1.94%  '
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('D8').Value = '''41.48'
$ws.Range('E8').Value = '  +5.07%  '
$ws.Range('E9').Value = '  +3.83%  '
$ws.Range('E10').Value = '  +2.07%  '
$ws.Range('D11').Value = '''0.0984'
$ws.Range('E11').Value = '  -1.06%  '
$ws.Range('D12').Value = '''2.112.76'
$ws.Range('E12').Value = '  +2.01%  '
$ws.Range('E13').Value = '  +5.06%  '
$ws.Range('D14').Value = '''1.848.28'
$ws.Range('E14').Value = '  +2.49%  '
$ws.Range('D15').Value = '''0.674'
$ws.Range('E15').Value = '  +1.99%  '
$ws.Range('D16').Value = '''4.67'
$ws.Range('E16').Value = '  +2.54%  '
$ws.Range('D17').Value = '''34.985.34'
$ws.Range('E17').Value = '  +0.22%  '
$ws.Range('D18').Value = '''70.07'
$ws.Range('E18').Value = '  +1.44%  '
$ws.Range('E19').Value = '  +1.61%  '
$ws.Range('D20').Value = '''240.75'
$ws.Range('E20').Value = '  +1.07%  '
$ws.Range('E21').Value = '  +3.64%  '
$ws.Range('E22').Value = '  +3.11%  '
$ws.Range('E23').Value = '  +0.17%  '
$ws.Range('E24').Value = '  +1.00%  '
$ws.Range('D25').Value = '''172.75'
$ws.Range('E25').Value = '  -0.28%  '
$ws.Range('D26').Value = '''7.83'
$ws.Range('E26').Value = '  +0.53%  '
$ws.Range('D27').Value = '''17.52'
$ws.Range('E27').Value = '  +1.87%  '
$ws.Range('D28').Value = '''0.124'
$ws.Range('E28').Value = '  +3.71%  '
$ws.Range('D29').Value = '''1.69'
$ws.Range('E29').Value = '  +9.63%  '
$ws.Range('E30').Value = '  +0.13%  '
$ws.Range('D31').Value = '''0.0554'
$ws.Range('E31').Value = '  +1.31%  '
$ws.Range('E32').Value = '  -0.03%  '
$ws.Range('E33').Value = '  +0.77%  '
$ws.Range('E34').Value = '  +25.12%  '
$ws.Range('D35').Value = '''1.95'
$ws.Range('E35').Value = '  +11.11%  '
$ws.Range('B36').Value = 'TrustWalletToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D36').Value = '''1.24'
$ws.Range('E36').Value = '  +6.79%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = '''0.747'
$ws.Range('E37').Value = '  +8.69%  '
$ws.Range('E38').Value = '  +11.65%  '
$ws.Range('D39').Value = '''90.04'
$ws.Range('E39').Value = '  -1.17%  '
$ws.Range('D40').Value = '''1.348.60'
$ws.Range('E40').Value = '  +3.07%  '
$ws.Range('E41').Value = '  +2.80%  '
$ws.Range('E42').Value = '  +3.30%  '
$ws.Range('D43').Value = '''2.28'
$ws.Range('E43').Value = '  +3.52%  '
$ws.Range('E44').Value = '  -1.64%  '
$ws.Range('E45').Value = '  +2.62%  '
$ws.Range('E46').Value = '  +4.22%  '
$ws.Range('D47').Value = '''6.34'
$ws.Range('E47').Value = '  +3.65%  '
$ws.Range('D48').Value = '''2.033.31'
$ws.Range('E48').Value = '  +2.07%  '
$ws.Range('E49').Value = '  +17.77%  '
$ws.Range('E50').Value = '  +0.18%  '
$ws.Range('E51').Value = '  -0.49%  '
